$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "Nose"

$ws.Range("H7").Value = "RBigToe"
$ws.Range("H8").Value = "RSmallToe"

$ws.Range("H13").Value = "LBigToe"
$ws.Range("H14").Value = "LSmallToe"

$ws.Range("E15").Value = "LShoulder"
$ws.Range("E16").Value = "RShoulder"
$ws.Range("H16").Value = "Neck"
$ws.Range("E17").Value = "LElbow"
$ws.Range("H17").Value = "Head"
$ws.Range("E18").Value = "RElbow"
$ws.Range("E19").Value = "LWrist"
$ws.Range("E20").Value = "RWrist"

$ws.Range("E27").Value = "LHip"
$ws.Range("E28").Value = "RHip"
$ws.Range("E29").Value = "LKnee"
$ws.Range("E30").Value = "RKnee"
$ws.Range("E31").Value = "LAnkle"
$ws.Range("E32").Value = "RAnkle"
$ws.Range("E33").Value = "LHeel"
$ws.Range("E34").Value = "RHeel"

$ws.Range("E17").Select()
